$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*Deixar o mesmo nome para Usuário / Funcionários*") {
        $p.Range.Font.Color = 1974729
    }
}
